# "Doing Updates for Financials" - refresh the yearly figures on the FONR sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FONR")

# Inventory (row 44): update yearly figures
$ws.Range("D44").Value = 1500
$ws.Range("E44").Value = 2400
$ws.Range("F44").Value = 2100
$ws.Range("G44").Value = 2900
$ws.Range("H44").Value = 3200
$ws.Range("I44").Value = 2500
$ws.Range("J44").Value = 3300

# Other Current Assets (row 45): update yearly figures
$ws.Range("D45").Value = 1300
$ws.Range("E45").Value = 1300
$ws.Range("F45").Value = 800
$ws.Range("G45").Value = 900
$ws.Range("H45").Value = 1000
$ws.Range("I45").Value = 1100
$ws.Range("J45").Value = 300

# Capital Expenditures (row 91): update yearly figures
$ws.Range("D91").Value = -2800
$ws.Range("E91").Value = -2900
$ws.Range("F91").Value = -700
$ws.Range("G91").Value = -100
$ws.Range("H91").Value = -600
$ws.Range("I91").Value = -1100
$ws.Range("J91").Value = -1100
